# The commit swaps the two embedded theme parts: the theme that is bound to
# the slide master (theme1.xml, originally "Integral") becomes the "Office
# Theme" palette that used to live on the notes-master theme (theme2.xml),
# and vice versa. The only thing that actually differs between the two
# theme parts in this deck is their <a:clrScheme> (the font scheme and
# format scheme are byte-identical already) plus some cosmetic name
# attributes, so re-pointing the slide master's 12 theme colours at the
# "Office Theme" palette reproduces the meaningful, visible part of that
# swap.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Target palette = the "Office Theme" colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) expressed as COM RGB ints (0xBBGGRR).
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72

# Best-effort: try to relabel the colour scheme / theme names to match the
# "Office Theme" naming too (some hosts expose these as writable; if not,
# this is a harmless no-op under try/catch).
try { $tcs.Name = "Office" } catch { }
try { $theme.Name = "Office Theme" } catch { }
